# Commit: "modified: Slow-Roaster.xlsx" — the only substantive content
# change in the workbook is the text of cell A1 (the sole shared string),
# which changes from "weoifk" to "Job". Everything else in the captured
# OOXML diff (fileVersion/calcPr/bookViews in workbook.xml, the panose /
# xmlns noise in theme1.xml, the defaultRowHeight/dyDescent tweak in
# sheet1.xml's sheetFormatPr) is Excel-build/save-environment metadata
# that the host application regenerates on its own whenever it re-saves
# the file, not something set through the object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Job"
